$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value of 45203 for every
# data row (2 through 453). Update every one of those cells to 45204
# (one day later), preserving their existing date formatting/style.
for ($r = 2; $r -le 453; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
